$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 198.33333
$ws.Range("I115").Value = 198.33333
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 594.99999
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 972.00001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7888.857
$ws.Range("I132").Value = 8164.3706
$ws.Range("J132").Value = 450
$ws.Range("K132").Value = 24493.1118
$ws.Range("L132").Value = 1350
$ws.Range("M132").Value = -21963.1118
$ws.Range("N132").Value = -6410

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3568.923
$ws.Range("I138").Value = 2931.5264
$ws.Range("J138").Value = 3935.9092
$ws.Range("K138").Value = 8794.5792
$ws.Range("L138").Value = 11807.7276
$ws.Range("M138").Value = -3654.5792
$ws.Range("N138").Value = -22087.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1238004.5
$ws.Range("I32").Value = 1323940
$ws.Range("J32").Value = 34908
$ws.Range("K32").Value = 1323940
$ws.Range("L32").Value = 34908
$ws.Range("M32").Value = -1323653
$ws.Range("N32").Value = -35482

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 753695.75
$ws.Range("I74").Value = 838243.1
$ws.Range("J74").Value = 20952.166
$ws.Range("K74").Value = 838243.1
$ws.Range("L74").Value = 20952.166
$ws.Range("M74").Value = -837369.1
$ws.Range("N74").Value = -22700.166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 753695.75
$ws.Range("I77").Value = 838243.1
$ws.Range("J77").Value = 20952.166
$ws.Range("K77").Value = 4191215.5
$ws.Range("L77").Value = 104760.83
$ws.Range("M77").Value = -4186847.5
$ws.Range("N77").Value = -113496.83

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 45000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 45000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 45000
$ws.Range("N114").Value = -53678

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 18718
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 18718
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 18718
$ws.Range("N125").Value = -28558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 50400
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 50400
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 50400
$ws.Range("N130").Value = -60440

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2442.1724
$ws.Range("I132").Value = 1928.8889
$ws.Range("J132").Value = 3282.0908
$ws.Range("K132").Value = 5786.6667
$ws.Range("L132").Value = 9846.2724
$ws.Range("M132").Value = -3256.6667
$ws.Range("N132").Value = -14906.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3206549.8
$ws.Range("I134").Value = 1499.0834
$ws.Range("J134").Value = 41667156
$ws.Range("K134").Value = 4497.2502
$ws.Range("L134").Value = 125001468
$ws.Range("M134").Value = -1962.2502
$ws.Range("N134").Value = -125006538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7782031
$ws.Range("I31").Value = 4447796
$ws.Range("J31").Value = 11116266
$ws.Range("K31").Value = 4447796
$ws.Range("L31").Value = 11116266
$ws.Range("M31").Value = -4447501
$ws.Range("N31").Value = -11116856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 2600
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 2600
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 2600
$ws.Range("N33").Value = -3358

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7782031
$ws.Range("I34").Value = 4447796
$ws.Range("J34").Value = 11116266
$ws.Range("K34").Value = 4447796
$ws.Range("L34").Value = 11116266
$ws.Range("M34").Value = -4447594
$ws.Range("N34").Value = -11116670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 35443.633
$ws.Range("I86").Value = 103553.625
$ws.Range("J86").Value = 10676.363
$ws.Range("K86").Value = 103553.625
$ws.Range("L86").Value = 10676.363
$ws.Range("M86").Value = -102430.625
$ws.Range("N86").Value = -12922.363

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 35443.633
$ws.Range("I89").Value = 103553.625
$ws.Range("J89").Value = 10676.363
$ws.Range("K89").Value = 517768.125
$ws.Range("L89").Value = 53381.815
$ws.Range("M89").Value = -512152.125
$ws.Range("N89").Value = -64613.815

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 16227.571
$ws.Range("I94").Value = 50999
$ws.Range("J94").Value = 2319
$ws.Range("K94").Value = 50999
$ws.Range("L94").Value = 2319
$ws.Range("M94").Value = -50548
$ws.Range("N94").Value = -3221

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2922.1904
$ws.Range("I132").Value = 2734.3157
$ws.Range("J132").Value = 4707
$ws.Range("K132").Value = 8202.947100000001
$ws.Range("L132").Value = 14121
$ws.Range("M132").Value = -5672.947100000001
$ws.Range("N132").Value = -19181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3652.7703
$ws.Range("I134").Value = 2463.4866
$ws.Range("J134").Value = 4842.054
$ws.Range("K134").Value = 7390.459800000001
$ws.Range("L134").Value = 14526.162
$ws.Range("M134").Value = -4855.459800000001
$ws.Range("N134").Value = -19596.162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1482.3334
$ws.Range("I80").Value = 1482.3334
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1482.3334
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -484.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1482.3334
$ws.Range("I83").Value = 1482.3334
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7411.666999999999
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -2419.666999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1449.6
$ws.Range("I113").Value = 1449.6
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1449.6
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 720.4000000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H129").Value = 55333.332
$ws.Range("I129").Value = 48000
$ws.Range("J129").Value = 59000
$ws.Range("K129").Value = 48000
$ws.Range("L129").Value = 59000
$ws.Range("M129").Value = -43000
$ws.Range("N129").Value = -69000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 80000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 80000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3631.6667
$ws.Range("I7").Value = 2872.5
$ws.Range("J7").Value = 5150
$ws.Range("K7").Value = 2872.5
$ws.Range("L7").Value = 5150
$ws.Range("M7").Value = -2760.5
$ws.Range("N7").Value = -5374

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3954.125
$ws.Range("I82").Value = 4199.25
$ws.Range("J82").Value = 3709
$ws.Range("K82").Value = 4199.25
$ws.Range("L82").Value = 3709
$ws.Range("M82").Value = -3838.25
$ws.Range("N82").Value = -4431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3954.125
$ws.Range("I85").Value = 4199.25
$ws.Range("J85").Value = 3709
$ws.Range("K85").Value = 4199.25
$ws.Range("L85").Value = 3709
$ws.Range("M85").Value = -2951.25
$ws.Range("N85").Value = -6205

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3631.6667
$ws.Range("I126").Value = 2872.5
$ws.Range("J126").Value = 5150
$ws.Range("K126").Value = 8617.5
$ws.Range("L126").Value = 15450
$ws.Range("M126").Value = -6147.5
$ws.Range("N126").Value = -20390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 9999.666999999999
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 9999.666999999999
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 9999.666999999999
$ws.Range("N29").Value = -10579.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8774660
$ws.Range("I132").Value = 11907353
$ws.Range("J132").Value = 3119.8
$ws.Range("K132").Value = 35722059
$ws.Range("L132").Value = 9359.400000000001
$ws.Range("M132").Value = -35719529
$ws.Range("N132").Value = -14419.4
